$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")
$rushing.Range("C2").Value = 2
$rushing.Range("D2").Value = 1
$rushing.Range("F2").Value = 3

$rushing.Range("C6").Value = 28
$rushing.Range("D6").Value = 22
$rushing.Range("F6").Value = 4

$rushing.Range("D7").Value = 16
$rushing.Range("F7").Value = 10

$rushing.Range("C10").Value = 7

$rushing.Range("C12").Value = 4

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

$receiving.Range("C2").Value = 14
$receiving.Range("D2").Value = 10
$receiving.Range("G2").Value = 5
$receiving.Range("H2").Value = 4

$receiving.Range("C3").Value = 39
$receiving.Range("D3").Value = 25

$receiving.Range("C5").Value = 3
$receiving.Range("D5").Value = 2

$receiving.Range("C6").Value = 42
$receiving.Range("D6").Value = 25

$receiving.Range("C7").Value = 74
$receiving.Range("D7").Value = 57
$receiving.Range("E7").Value = 9
$receiving.Range("G7").Value = 12
$receiving.Range("H7").Value = 9

$receiving.Range("C9").Value = 54
$receiving.Range("D9").Value = 33
$receiving.Range("E9").Value = 21
$receiving.Range("F9").Value = 8
$receiving.Range("G9").Value = 7
$receiving.Range("H9").Value = 4

$receiving.Range("C10").Value = 29
$receiving.Range("D10").Value = 21

$receiving.Range("C11").Value = 9
$receiving.Range("D11").Value = 4

$receiving.Range("C15").Value = 33
$receiving.Range("D15").Value = 21
$receiving.Range("G15").Value = 6
$receiving.Range("H15").Value = 3
